# Upload Leave Card 12/27/2023 4:01 PM
# - Fill in EARNED (VL) for Sep/Oct/Nov 2023 rows (67-69)
# - Fill in the Dec 2023 row (70) particulars/undertime/remarks for a forced-leave entry
# - Insert a new "2024" year-marker row (71), pushing the existing 2024/2025 rows down
# - Keep the Table15 table range and trailing total-row formatting in sync

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Insert a new row for the "2024" year marker, shifting everything below down ---
$ws.Rows.Item(71).Insert()

# Keep the structured table in sync with the newly-inserted row
$lo = $ws.ListObjects.Item("Table15")
$lo.Resize($ws.Range("A8:K118"))

# Copy the formatting from the existing "2023" year-marker row (58) onto the new row
$ws.Range("A58:K58").Copy() | Out-Null
$ws.Range("A71:K71").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the year-marker label text and restore the calculated-column formula
$ws.Range("A71").Value = "'2024"
$ws.Range("G71").Formula = '=IF(ISBLANK(Table15[[#This Row],[EARNED]]),"",Table15[[#This Row],[EARNED]])'

# Restore the calculated-column formula on the new trailing total row created by Resize
$ws.Range("G118").Formula = '=IF(ISBLANK(Table15[[#This Row],[EARNED]]),"",Table15[[#This Row],[EARNED]])'

# --- Fill in the December 2023 (row 70) leave entry ---
$ws.Range("B70").Value = "FL(2-0-0)"
$ws.Range("D70").Value = 2
$ws.Range("K70").Value = "12/12,13/2023"

# --- Fill EARNED for the three already-elapsed 2023 months ---
$ws.Range("C67").Value = 1.25
$ws.Range("C68").Value = 1.25
$ws.Range("C69").Value = 1.25

# Match the final scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("G75").Select() | Out-Null
